# Re-pull / push updated "dSF" (column F) data for each row on the active sheet.
# Commit message: "repull data, push all data, mean calculation"
# The diff shows column F ("dSF") values being refreshed (repulled) with
# newly computed deltas for most data rows, while column E ("dS0") is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -7
    4  = -3
    5  = 2
    6  = -4
    7  = -1
    8  = 3
    9  = -8
    10 = -6
    11 = 1
    12 = -1
    14 = -1
    15 = 5
    16 = -6
    17 = 6
    19 = -1
    20 = 2
    21 = 5
    22 = 2
    23 = 1
    24 = 3
    25 = 3
    26 = -1
    27 = -1
    28 = -3
    29 = 4
    30 = 4
    31 = -1
    32 = -3
    33 = 2
    34 = -3
    35 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
